$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the example filenames in row 2
$ws.Range("B2").Value = "niches.tsv"
$ws.Range("C2").Value = "digital_humanities.csv"

# Update the active selection to match the target file (C2)
$ws.Range("C2").Select()
